$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180, pushing the existing rows 180-209 down to 181-210.
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new data record.
$ws.Cells.Item(180, 1).Value = 10
$ws.Cells.Item(180, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(180, 3).Value = "La Araucanía"
$ws.Cells.Item(180, 4).Value = 44522
$ws.Cells.Item(180, 5).Value = 9
$ws.Cells.Item(180, 6).Value = 100112017
$ws.Cells.Item(180, 7).Value = "Apio"
$ws.Cells.Item(180, 8).Value = "Americana (o)"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 30
$ws.Cells.Item(180, 11).Value = 9000
$ws.Cells.Item(180, 12).Value = 9000
$ws.Cells.Item(180, 13).Value = 9000
$ws.Cells.Item(180, 14).Value = "$/docena de matas"
$ws.Cells.Item(180, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(180, 16).Value = 1500
$ws.Cells.Item(180, 17).Value = 6
$ws.Cells.Item(180, 18).Value = "Hortaliza"
